$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IFrames_WYSIWYGEditor")

# Update the Special_Effects value for TC2 (row 3, column D)
$ws.Range("D3").Value = "Bold;Justify;right;Decrease-indent;Italic"

# Update the active cell selection to A3:D3
$ws.Range("A3:D3").Select()
